{"js": "// Highlight the \"MVC\" (SNO 5) and \"NodeJS\" (SNO 16) rows of the syllabus\n// table in green (RGB 00B050), matching the existing \"Two Way Data\n// Binding\" row's styling. Only the first three cells of each row\n// (the number, the short topic name, and the full topic description)\n// carry this formatting in the source document; the merged duration\n// cell on the right is left untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row indexes (0-based) of the rows whose topic text reads \"5 / MVC /\n// MVC in AngularJS\" and \"16 / NodeJS / Interacting With NodeJS\".\nconst targetRowIndexes = [8, 19];\nconst targetColumnIndexes = [0, 1, 2];\n\nfor (const rowIndex of targetRowIndexes) {\n  for (const columnIndex of targetColumnIndexes) {\n    const cell = table.getCellOrNullObject(rowIndex, columnIndex);\n    cell.body.font.color = \"#00B050\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight the \"MVC\" (SNO 5) and \"NodeJS\" (SNO 16) rows of the syllabus\n# table in green (RGB 00B050 == 0x00B050), matching the existing\n# \"Two Way Data Binding\" row's styling. Only the first three cells of\n# each row (the number, the short topic name, and the full topic\n# description) carry this formatting in the source document; the merged\n# duration cell on the right is left untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word's Font.Color expects a BGR-packed long: R + G*256 + B*65536.\n# 0x00, 0xB0, 0x50 -> 0x0050B000 == 5287936 (hex string \"00B050\").\n$green = 5287936\n\n# Table.Cell(row, column) is 1-based: row 9 = SNO \"5\" / \"MVC\" / \"MVC in\n# AngularJS\"; row 20 = SNO \"16\" / \"NodeJS\" / \"Interacting With NodeJS\".\n$targetRows = @(9, 20)\n$targetColumns = @(1, 2, 3)\n\nforeach ($rowIndex in $targetRows) {\n    foreach ($columnIndex in $targetColumns) {\n        $cell = $t.Cell($rowIndex, $columnIndex)\n        $cell.Range.Font.Color = $green\n    }\n}\n"}
